$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data refresh.
# NumberFormat "@" + Style "Normal" keeps values as text (matching original inlineStr cells)
# while avoiding a stray number-format style being attached to the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.010.91'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.60%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.812.83'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +6.56%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '427.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +7.87%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.36'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.12%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.810.33'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +6.78%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.87%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.03%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -7.65%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.67%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.50'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +12.07%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.412.35'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +7.10%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +18.70%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.22%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.823.01'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +5.80%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '20.06'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.50%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +7.69%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '66.194.99'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.75%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '414.79'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.64%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.28'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +8.27%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.68'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.00%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +7.69%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '37.09'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +7.37%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.06'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +12.67%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +8.25%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.41'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.37%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +34.50%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '14.18'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +17.83%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '710.35'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.29%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +12.23%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.76'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.17%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +37.55%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '38.82'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.30%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '55.97'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.20%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0468'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +5.44%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₃0723'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +12.86%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.88'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.85%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.82%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.26'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +6.20%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +8.56%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.322'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +15.17%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.40'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +37.16%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.51%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.82%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '141.60'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.28%  '
$ws.Range('E51').Style = 'Normal'
